$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Registro_acciones_inventario": append a new row (row 27) of data
# ---------------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Registro_acciones_inventario")

$wsReg.Cells.Item(27, 1).Value = 46
$wsReg.Cells.Item(27, 2).Value = 45836
$wsReg.Cells.Item(27, 2).NumberFormat = "YYYY-MM-DD"
$wsReg.Cells.Item(27, 3).Value = "zeta"
$wsReg.Cells.Item(27, 4).Value = "Hornos"
$wsReg.Cells.Item(27, 5).Value = 2
$wsReg.Cells.Item(27, 6).Value = 2
$wsReg.Cells.Item(27, 7).Value = 4
$wsReg.Cells.Item(27, 8).Value = "Activo"
$wsReg.Cells.Item(27, 9).Value = "Soluciones Electromecánicas"
$wsReg.Cells.Item(27, 10).Value = "Almacén 4"
$wsReg.Cells.Item(27, 11).Value = "Inventario 2"
$wsReg.Cells.Item(27, 12).Value = "REGISTRAR"
$wsReg.Cells.Item(27, 13).Value = "nestor"

# ---------------------------------------------------------------------------
# Sheet "Usuarios": update connection/permission bookkeeping fields
# ---------------------------------------------------------------------------
$wsUsr = $wb.Worksheets.Item("Usuarios")

# nestor (row 2) - refreshed last-connection timestamp
$wsUsr.Cells.Item(2, 7).Value = 45838.59877830116

# piero (row 3) - refreshed last-connection timestamp, permissions revoked
$wsUsr.Cells.Item(3, 7).Value = 45837.00596537504
$wsUsr.Cells.Item(3, 8).Value = $false
$wsUsr.Cells.Item(3, 9).Value = $false
$wsUsr.Cells.Item(3, 10).Value = $false

# luis (row 4) - refreshed last-connection timestamp
$wsUsr.Cells.Item(4, 7).Value = 45837.00615547206

# johnny (row 6) - permissions granted
$wsUsr.Cells.Item(6, 8).Value = $true
$wsUsr.Cells.Item(6, 9).Value = $true
$wsUsr.Cells.Item(6, 10).Value = $true
